$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.433526277542114
$ws.Range("B1").Value = 1.222121238708496
$ws.Range("C1").Value = 5.245344638824463
$ws.Range("D1").Value = 3.587836027145386
$ws.Range("E1").Value = 0.6805036067962646
